$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16 (shifts existing rows 16-30 down to 17-31),
# inheriting formatting/styles from the row above (row 15), which matches
# the target (H/I columns keep style index 1, others unstyled).
$ws.Rows.Item(16).Insert() | Out-Null

# Fill in the new BOM entry for the missing ceramic capacitor.
# Values are set in the same order the target sharedStrings table grows
# (A, B, H, I, then E) so new unique strings are appended in the right order.
$ws.Cells.Item(16, 1).Value = "C601,C602,C701,C702"
$ws.Cells.Item(16, 2).Value = "GRM18R60J105KA01J"
$ws.Cells.Item(16, 3).Value = 4
$ws.Cells.Item(16, 8).Value = "https://www.digikey.de/product-detail/de/murata-electronics-north-america/GRM188R60J105KA01J/490-6404-1-ND/3845601"
$ws.Cells.Item(16, 9).Value = "https://www.mouser.de/ProductDetail/Murata/GRM188R60J105KA01J/"
$ws.Cells.Item(16, 5).Value = "1uF / 6.3V"
$ws.Cells.Item(16, 4).Value = "Ceramic capacitor, X5R"
$ws.Cells.Item(16, 6).Value = "SMD0603"

# The row-insert does not automatically shift the existing Hyperlinks
# collection's ranges, so rebuild it from scratch with the correct
# (possibly shifted) target cells. Rows below the inserted row (>=16 in the
# pre-insert numbering) move down by one; rows above stay the same.
$ws.Hyperlinks.Delete() | Out-Null

$ws.Hyperlinks.Add($ws.Range("H2"), "http://www.digikey.de/scripts/DkSearch/dksus.dll?Detail&itemSeq=234870402", [Type]::Missing, [Type]::Missing, "http://www.digikey.de/scripts/DkSearch/dksus.dll?Detail&itemSeq=234870402") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "http://www.mouser.de/ProductDetail/Intel/10CL010YE144C8G/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H3"), "https://www.digikey.de/product-detail/de/analog-devices-inc/ADV7125KSTZ140/ADV7125KSTZ140-ND/654207") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "http://www.mouser.de/ProductDetail/Analog-Devices/ADV7125KSTZ140/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I4"), "http://www.mouser.de/ProductDetail/STMicroelectronics/TSH122ICT/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H4"), "https://www.digikey.de/product-detail/de/stmicroelectronics/TSH122ICT/497-8332-1-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I5"), "http://www.mouser.de/ProductDetail/Texas-Instruments/SN74LVC3G17DCTR/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I10"), "http://www.mouser.de/ProductDetail/Murata/GRM219R60J476ME44D/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H10"), "https://www.digikey.de/product-detail/de/murata-electronics-north-america/GRM219R60J476ME44D/490-13249-1-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H12"), "https://www.digikey.de/product-detail/de/murata-electronics-north-america/GRM21BR60J106ME19L/490-1718-1-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I12"), "http://www.mouser.de/ProductDetail/Murata/GRM21BR60J106ME19L/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I6"), "http://www.mouser.de/ProductDetail/Cypress/S25FL116K0XMFB013/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H7"), "https://www.digikey.de/product-detail/de/texas-instruments/TLV70012DDCR/296-39275-1-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I7"), "http://www.mouser.de/ProductDetail/Texas-Instruments/TLV70012DDCR/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H8"), "https://www.digikey.de/product-detail/de/texas-instruments/TLV70025DDCR/296-32411-1-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I8"), "http://www.mouser.de/ProductDetail/Texas-Instruments/TLV70025DDCR/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H9"), "https://www.digikey.de/product-detail/de/abracon-llc/ASE-50.000MHZ-LC-T/535-9577-1-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I9"), "http://www.mouser.de/ProductDetail/ABRACON/ASE-50000MHZ-L-C-T/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H11"), "https://www.digikey.de/product-detail/de/murata-electronics-north-america/GRM188R71H104KA93J/490-9735-1-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I11"), "http://www.mouser.de/ProductDetail/Murata/GRM188R71H104KA93J/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I13"), "http://www.mouser.de/ProductDetail/Murata/GRM188R70J103KA01D/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H13"), "https://www.digikey.de/product-detail/de/murata-electronics-north-america/GRM188R70J103KA01D/490-9729-1-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I17"), "http://www.mouser.de/ProductDetail/TDK/MPZ1608S221A/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H14"), "https://www.digikey.de/product-detail/de/murata-electronics-north-america/GRM32DR60J336ME19L/490-3389-1-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I14"), "http://www.mouser.de/ProductDetail/Murata-Electronics/GRM32DR60J336ME19L/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I15"), "http://www.mouser.de/ProductDetail/Murata/GRM31CR60J226KE19L/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H15"), "https://www.digikey.de/product-detail/de/murata-electronics-north-america/GRM31CR60J226KE19L/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H20"), "https://www.digikey.de/product-detail/de/yageo/RC0603FR-0710KL/311-10.0KHRCT-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I20"), "http://www.mouser.de/ProductDetail/Yageo/RC0603FR-0710KL/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H19"), "https://www.digikey.de/product-detail/de/yageo/RC0603FR-071KL/311-1.00KHRCT-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I19"), "http://www.mouser.de/ProductDetail/Yageo/RC0603FR-071KL/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H22"), "https://www.digikey.de/product-detail/de/yageo/RC0603FR-074K7L/311-4.70KHRCT-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I22"), "http://www.mouser.de/ProductDetail/Yageo/RC0603FR-074K7L/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I24"), "http://www.mouser.de/ProductDetail/Yageo/RC0603FR-07475RL/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H24"), "https://www.digikey.de/product-detail/de/yageo/RC0603FR-07475RL/311-475HRCT-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I23"), "http://www.mouser.de/ProductDetail/Yageo/RC0603FR-0775RL/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H23"), "https://www.digikey.de/product-detail/de/yageo/RC0603FR-0775RL/311-75.0HRCT-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H26"), "https://www.digikey.de/product-detail/de/bourns-inc/CAT16-47R0F4LF/CAT16-47R0F4LFCT-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I26"), "http://www.mouser.de/ProductDetail/Bourns/CAT16-47R0F4LF/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I18"), "http://www.mouser.de/ProductDetail/3M/30310-6002HB/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H18"), "https://www.digikey.de/product-detail/de/3m/30310-6002HB/30310-6002HB-ND/1237393") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H21"), "https://www.digikey.de/product-detail/de/yageo/RC0603FR-07536RL/311-536HRCT-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I25"), "http://www.mouser.de/ProductDetail/Yageo/RC0603FR-0724R9L/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H25"), "https://www.digikey.de/product-detail/de/yageo/RC0603FR-0724R9L/311-24.9HRCT-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I27"), "http://www.mouser.de/search/ProductDetail.aspx?R=0virtualkey0virtualkeyCAT16-75R0F4LF") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H27"), "https://www.digikey.de/product-detail/de/bourns-inc/CAT16-750J4LF/CAT16-750J4LFCT-ND/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I21"), "http://www.mouser.de/ProductDetail/Yageo/RC0603FR-07536RL/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H6"), "https://www.digikey.de/product-detail/de/cypress-semiconductor-corp/S25FL116K0XMFA043/428-4061-1-ND/") | Out-Null

# Match the selection recorded in the target sheet view.
$ws.Range("E17").Select() | Out-Null
